$d = $word.ActiveDocument

# Helper-free approach: for each target paragraph, run a *scoped* Find/Replace
# (Find executed on the paragraph's own Range) so formatting is inherited
# from the matched run and so duplicate phrases elsewhere in the document
# are not touched.

# Paragraph indices (1-based, from $d.Paragraphs) as found in before.docx:
# 45 Definir los requerimientos funcionales
# 46 Identificacion los requerimientos funcionales a automatizar
# 47 Eleccion de recursos tecnicos (Definir los requerimientos no funcionale).
# 48 Aprendizaje de los recursos tecnicos  (trailing space)
# 50 Realizar el modelo E-R
# 51 Definir validaciones y funcionalidades
# 52 Diagramar los procesos de negocio a automatizar
# 53 Realizar el modelo logico
# 54 Elegir el estandar de usabilidad
# 56 Implementar la base de datos
# 57 Insertar informacion en la base de datos
# 58 Implementar los modulos de mantenimiento de la informacion
# 59 Implementar los modulos de la automatizacion del proceso de negocio
# 60 Implementacion de los reportes
# 61 Disenar las interfaces de usuario
# 64 bibliography paragraph

# 1) Definir los requerimientos funcionales -> append description
$d.Paragraphs(45).Range.Find.Execute("Definir los requerimientos funcionales", $true, $false, $false, $false, $false, $true, 1, $false, "Definir los requerimientos funcionales: Extraer de los procesos de negocio las funcionalidades a ser implementadas.", 2)

# 2) Identificacion ... automatizar -> append description
$d.Paragraphs(46).Range.Find.Execute(" automatizar", $true, $false, $false, $false, $false, $true, 1, $false, " automatizar: De las funcionalidades obtenida filtrar las que se pueden implementar a nivel de software.", 2)

# 3) Eleccion de recursos tecnicos (...).  -> replace trailing "." with ": Analizar...requerimientos."
$d.Paragraphs(47).Range.Find.Execute("no funciónale).", $true, $false, $false, $false, $false, $true, 1, $false, "no funciónale): Analizar en base a los requerimientos funcionales, cual seria el manejador mas optimo y el lenguaje de programación mas adecuado para implementar dichos requerimientos.", 2)

# 4) Aprendizaje de los recursos tecnicos (trailing space) -> insert description text before the trailing space
$d.Paragraphs(48).Range.Find.Execute("Aprendizaje de los recursos técnicos ", $true, $false, $false, $false, $false, $true, 1, $false, "Aprendizaje de los recursos técnicos: Periodo de tiempo necesario para familiarizarse con las herramientas necesarias para el desarrollo del software. ", 2)

# 5) Realizar el modelo E-R -> append description
$d.Paragraphs(50).Range.Find.Execute("odelo E-R", $true, $false, $false, $false, $false, $true, 1, $false, "odelo E-R: Basado en los requerimientos funcionales, se elabora un modelo que describe la estructura de la base de datos.", 2)

# 6) Definir validaciones y funcionalidades -> append description
$d.Paragraphs(51).Range.Find.Execute("Definir validaciones y funcionalidades", $true, $false, $false, $false, $false, $true, 1, $false, "Definir validaciones y funcionalidades: Definir funcionalidades que no pueden ser modelados en el E-R.", 2)

# 7) Diagramar los procesos de negocio a automatizar -> remove jc=both, append description
$p52 = $d.Paragraphs(52)
$p52.Alignment = 0
$p52.Range.Find.Execute("Diagramar los procesos de negocio a automatizar", $true, $false, $false, $false, $false, $true, 1, $false, "Diagramar los procesos de negocio a automatizar: Realizar diagramas que definan la lógica y los procesos necesarios definir funcionalidades no definidas por el E-R.", 2)

# 8) Realizar el modelo logico -> append description
$d.Paragraphs(53).Range.Find.Execute("Realizar el modelo lógico", $true, $false, $false, $false, $false, $true, 1, $false, "Realizar el modelo lógico: Basado en el modelo E-R, se realiza un modelo mas riguroso y enfocado a la implementación de la base de datos.", 2)

# 9) Elegir el estandar de usabilidad -> append description
$d.Paragraphs(54).Range.Find.Execute("Elegir el estándar de usabilidad", $true, $false, $false, $false, $false, $true, 1, $false, "Elegir el estándar de usabilidad: Elegir un estándar de usabilidad que garantice la integridad de la ampliación y la comunicación efectiva del usuario.", 2)

# 10) Implementar la base de datos -> append description
$d.Paragraphs(56).Range.Find.Execute("Implementar la base de datos", $true, $false, $false, $false, $false, $true, 1, $false, "Implementar la base de datos: Basado en el modelo lógico se generan las tablas y relaciones.", 2)

# 11) Insertar información en la base de datos -> append description
$d.Paragraphs(57).Range.Find.Execute("Insertar información en la base de datos", $true, $false, $false, $false, $false, $true, 1, $false, "Insertar información en la base de datos: Mediante código SQL se inserta información en la tabla de la base de datos.", 2)

# 12) Paragraph 58 (was "Implementar los módulos de mantenimiento de la información")
#     -> becomes "Diseñar las interfaces de usuario: Usando los recursos técnicos se implementa la interfaz en base a los reportes necesarios y los procesos de negocio"
$d.Paragraphs(58).Range.Find.Execute("Implementar los módulos de mantenimiento de la información", $true, $false, $false, $false, $false, $true, 1, $false, "Diseñar las interfaces de usuario: Usando los recursos técnicos se implementa la interfaz en base a los reportes necesarios y los procesos de negocio", 2)

# 13) Paragraph 59 (was "I" + "mplementar los módulos de la automatización del proceso de negocio")
#     -> becomes "Implementar los módulos de mantenimiento de la información: Crear CRUD por cada tabla, asegundo que funcionan con sus validaciones respectivas."
$d.Paragraphs(59).Range.Find.Execute("Implementar los módulos de la automatización del proceso de negocio", $true, $false, $false, $false, $false, $true, 1, $false, "Implementar los módulos de mantenimiento de la información: Crear CRUD por cada tabla, asegundo que funcionan con sus validaciones respectivas.", 2)

# 14) Paragraph 60 (was "Implementación de los reportes")
#     -> becomes "Implementar los módulos de la automatización del proceso de negocio: Basado en los procesos de negocio, se crean los módulos necesarios."
$d.Paragraphs(60).Range.Find.Execute("Implementación de los reportes", $true, $false, $false, $false, $false, $true, 1, $false, "Implementar los módulos de la automatización del proceso de negocio: Basado en los procesos de negocio, se crean los módulos necesarios.", 2)

# 15) Paragraph 61 (was "Diseñar las interfaces de usuario")
#     -> becomes "Implementación de los reportes: Basado en procesos de negocio, se realizan las consultas necesarias a la base de datos y se muestran en la interfaz."
$d.Paragraphs(61).Range.Find.Execute("Diseñar las interfaces de usuario", $true, $false, $false, $false, $false, $true, 1, $false, "Implementación de los reportes: Basado en procesos de negocio, se realizan las consultas necesarias a la base de datos y se muestran en la interfaz.", 2)

Write-Output "done batch 1"
